# CRCP Sediment Data.xlsx - "Input data from Sept and new plots of raw data"
#
# This script reproduces the view/selection-state changes captured in the
# target diff:
#   - April2014LOI sheet: the frozen "bottomRight" pane's active cell moves
#     from Z20 to H17 (reviewing September data lower in the sheet).
#   - CRCP 2014 sheet: the sheet view no longer scrolls to topLeftCell E1,
#     and the selection moves from the single cell K9 to the C2:D20 block
#     (the newly populated SedPod/Tube raw-data range for the new plots).
#
# The CRCP 2014 sheet is left as the active/selected sheet at the end, to
# match its tabSelected="1" sheetView state in the workbook.

$wb = $excel.ActiveWorkbook

# --- April2014LOI: update the frozen pane's active selection ---------------
$loiSheet = $wb.Worksheets.Item("April2014LOI")
$loiSheet.Range("H17").Select()

# --- CRCP 2014: clear the scrolled topLeftCell and reselect C2:D20 --------
$crcpSheet = $wb.Worksheets.Item("CRCP 2014")
$crcpSheet.Range("C2:D20").Select()
